# Insert a new price-record row for "Apio" (Terminal Hortofrutícola Agro
# Chillán) right before the current row 86. This pushes the existing rows
# 86..213 down to 87..214 (dimension grows from A1:R213 to A1:R214) and
# the freshly inserted row 86 gets a brand-new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 7
$ws.Cells.Item(86, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(86, 3).Value = "Ñuble"
$ws.Cells.Item(86, 4).Value = 44757
$ws.Cells.Item(86, 5).Value = 16
$ws.Cells.Item(86, 6).Value = 100112017
$ws.Cells.Item(86, 7).Value = "Apio"
$ws.Cells.Item(86, 8).Value = "Americana (o)"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 160
$ws.Cells.Item(86, 11).Value = 8500
$ws.Cells.Item(86, 12).Value = 9000
$ws.Cells.Item(86, 13).Value = 8750
$ws.Cells.Item(86, 14).Value = "`$/docena de matas"
$ws.Cells.Item(86, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(86, 16).Value = 1458
$ws.Cells.Item(86, 17).Value = 6
$ws.Cells.Item(86, 18).Value = "Hortaliza"
